# Update the "想去人数" (want-to-go count) figures that were refreshed
# at data-generation time (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first 4 events
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 3205
$wsExpo.Range("F4").Value = 151
$wsExpo.Range("F5").Value = 22
$wsExpo.Range("F6").Value = 136

# Sheet "全部类型" (All types) - same events duplicated in the combined view
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3205
$wsAll.Range("F8").Value = 151
$wsAll.Range("F9").Value = 22
$wsAll.Range("F11").Value = 136
